# Commit: "Upload validation and error handling (#6327)"
#
# spec/support/roDuplicateDates.xlsx — third worksheet ("RO & CO Hearing
# Allocation") drops the "Central Office" allocation row from its data
# table and is renamed/retitled to drop the "Central Office" wording,
# since Central Office hearings are no longer tracked on this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# The data table's first data row (row 4) was the "Central Office" entry
# (City/State column blank, BFREGOFF column blank, just a day count).
# Remove it entirely so the Regional Office rows that followed (rows
# 5-61) shift up one row, and the trailing blank row (100) disappears
# along with it.
$ws.Rows.Item(4).Delete()

# Update the page title in A1 to drop the "and Central Office Hearings"
# wording, and rename the sheet tab to match.
$ws.Range("A1").Value = "Allocation of Regional Office Video Hearings"
$ws.Name = "RO Allocations"
